$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure price column (D) values are written as text to preserve exact formatting
foreach ($addr in @("D2", "D3", "D5", "D6", "D8", "D9", "D10", "D11", "D12", "D13", "D14", "D15", "D16", "D17", "D18", "D19", "D20", "D21", "D22", "D23", "D25", "D40", "D41", "D44", "D45", "D47", "D49", "D50")) {
    $ws.Range($addr).NumberFormat = "@"
}

# Ensure hour column (G) values are written as text to preserve exact formatting
$ws.Range("G2:G51").NumberFormat = "@"

# Apply cell value updates per diff
$ws.Range("D2").Value = "244.26"
$ws.Range("G2").Value = "14"
$ws.Range("D3").Value = "23.08"
$ws.Range("G3").Value = "14"
$ws.Range("G4").Value = "14"
$ws.Range("D5").Value = "0.05938"
$ws.Range("G5").Value = "14"
$ws.Range("D6").Value = "3.387"
$ws.Range("G6").Value = "14"
$ws.Range("G7").Value = "14"
$ws.Range("D8").Value = "0.9255"
$ws.Range("G8").Value = "14"
$ws.Range("B9").Value = "WazirX"
$ws.Range("C9").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D9").Value = "0.1413"
$ws.Range("E9").Value = "8WazirXWRX"
$ws.Range("G9").Value = "14"
$ws.Range("B10").Value = "MandalaExchangeToken"
$ws.Range("C10").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D10").Value = "0.07402"
$ws.Range("E10").Value = "9MandalaExchangeTokenMDX"
$ws.Range("G10").Value = "14"
$ws.Range("B11").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C11").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D11").Value = "0.03415"
$ws.Range("E11").Value = "10LiechtensteinCryptoassetsExchangeLCX"
$ws.Range("G11").Value = "14"
$ws.Range("B12").Value = "BitrueCoin"
$ws.Range("C12").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D12").Value = "0.03080"
$ws.Range("E12").Value = "11BitrueCoinBTR"
$ws.Range("G12").Value = "14"
$ws.Range("B13").Value = "BitMartToken"
$ws.Range("C13").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D13").Value = "0.09343"
$ws.Range("E13").Value = "12BitMartTokenBMX"
$ws.Range("G13").Value = "14"
$ws.Range("B14").Value = "MCDex"
$ws.Range("C14").Value = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
$ws.Range("D14").Value = "3.954"
$ws.Range("E14").Value = "13MCDexMCB"
$ws.Range("G14").Value = "14"
$ws.Range("B15").Value = "BitForexToken"
$ws.Range("C15").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D15").Value = "0.001590"
$ws.Range("E15").Value = "14BitForexTokenBF"
$ws.Range("G15").Value = "14"
$ws.Range("B16").Value = "CoinExToken"
$ws.Range("C16").Value = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
$ws.Range("D16").Value = "0.04807"
$ws.Range("E16").Value = "15CoinExTokenCET"
$ws.Range("G16").Value = "14"
$ws.Range("B17").Value = "One"
$ws.Range("C17").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Range("D17").Value = "0.0005944"
$ws.Range("E17").Value = "16OneONE"
$ws.Range("G17").Value = "14"
$ws.Range("D18").Value = "0.005509"
$ws.Range("G18").Value = "14"
$ws.Range("D19").Value = "0.004327"
$ws.Range("G19").Value = "14"
$ws.Range("D20").Value = "0.0009855"
$ws.Range("G20").Value = "14"
$ws.Range("D21").Value = "0.00007705"
$ws.Range("G21").Value = "14"
$ws.Range("D22").Value = "3.667"
$ws.Range("G22").Value = "14"
$ws.Range("D23").Value = "6.443"
$ws.Range("G23").Value = "14"
$ws.Range("G24").Value = "14"
$ws.Range("D25").Value = "0.3241"
$ws.Range("G25").Value = "14"
$ws.Range("G26").Value = "14"
$ws.Range("G27").Value = "14"
$ws.Range("G28").Value = "14"
$ws.Range("G29").Value = "14"
$ws.Range("G30").Value = "14"
$ws.Range("G31").Value = "14"
$ws.Range("G32").Value = "14"
$ws.Range("G33").Value = "14"
$ws.Range("G34").Value = "14"
$ws.Range("G35").Value = "14"
$ws.Range("G36").Value = "14"
$ws.Range("G37").Value = "14"
$ws.Range("G38").Value = "14"
$ws.Range("G39").Value = "14"
$ws.Range("D40").Value = "0.03920"
$ws.Range("G40").Value = "14"
$ws.Range("D41").Value = "0.006184"
$ws.Range("G41").Value = "14"
$ws.Range("G42").Value = "14"
$ws.Range("G43").Value = "14"
$ws.Range("D44").Value = "0.007298"
$ws.Range("E44").Value = "43LocalTradersLCTBestin24h"
$ws.Range("G44").Value = "14"
$ws.Range("D45").Value = "0.00005177"
$ws.Range("G45").Value = "14"
$ws.Range("G46").Value = "14"
$ws.Range("D47").Value = "0.0005804"
$ws.Range("G47").Value = "14"
$ws.Range("G48").Value = "14"
$ws.Range("D49").Value = "0.002312"
$ws.Range("G49").Value = "14"
$ws.Range("D50").Value = "0.00002101"
$ws.Range("G50").Value = "14"
$ws.Range("G51").Value = "14"
